# Insert a new weekly price-report row just above the current row 81
# (pushing the existing rows 81-97 down to 82-98), then populate the
# newly inserted row with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("81:81").Insert()

$ws.Cells.Item(81, 1).Value  = 10
$ws.Cells.Item(81, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(81, 3).Value  = "La Araucanía"
$ws.Cells.Item(81, 4).Value  = 44511
$ws.Cells.Item(81, 5).Value  = 9
$ws.Cells.Item(81, 6).Value  = 100112012
$ws.Cells.Item(81, 7).Value  = "Espinaca"
$ws.Cells.Item(81, 8).Value  = "Sin especificar"
$ws.Cells.Item(81, 9).Value  = "Primera"
$ws.Cells.Item(81, 10).Value = 50
$ws.Cells.Item(81, 11).Value = 8000
$ws.Cells.Item(81, 12).Value = 8000
$ws.Cells.Item(81, 13).Value = 8000
$ws.Cells.Item(81, 14).Value = "$/docena de atados"
$ws.Cells.Item(81, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(81, 16).Value = 2667
$ws.Cells.Item(81, 17).Value = 3
$ws.Cells.Item(81, 18).Value = "Hortaliza"
